# Gold standard mapping comparison
# Add a new "Final_Matches" worksheet at the end of the workbook summarizing
# the AZ/ASCTB label matches used for the gold-standard comparison.

$wb = $excel.ActiveWorkbook

# Grab the header formatting already used by the other sheets (bold, centered,
# bordered) so the new sheet's header matches without creating new styles.
$formatSource = $wb.Worksheets.Item(7).Range("A1:B1")

# Insert the new sheet after the last existing worksheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Final_Matches"

# Header row.
$ws.Range("A1").Value = "AZ.CT/LABEL"
$ws.Range("B1").Value = "ASCTB.CT/LABEL"

$formatSource.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Data rows: perfect label matches between AZ and ASCTB cell type labels.
$data = @(
    @("type B pancreatic cell", "type B pancreatic cell"),
    @("pancreatic A cell", "pancreatic A cell"),
    @("pancreatic D cell", "pancreatic D cell"),
    @("pancreatic acinar cell", "pancreatic acinar cell"),
    @("pancreatic ductal cell", "pancreatic ductal cell"),
    @("pancreatic PP cell", "pancreatic PP cell"),
    @("pancreatic stellate cell", "pancreatic stellate cell"),
    @("pancreatic epsilon cell", "pancreatic epsilon cell"),
    @("pancreatic endocrine cell", "pancreatic endocrine cell"),
    @("leukocyte", "lymphocyte"),
    @("leukocyte", "macrophage"),
    @("leukocyte", "dendritic cell"),
    @("leukocyte", "B cell"),
    @("leukocyte", "T cell")
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
